$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range's last row (data occupies columns A:D)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
